# Auto-generated edit script: applies cryptos.xlsx price/volume/row updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.705.52"
$ws.Range("E2").Value = "  +7.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.620.97"
$ws.Range("E3").Value = "  +7.42%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.87"
$ws.Range("E5").Value = "  +14.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "581.27"
$ws.Range("E6").Value = "  +3.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +4.18%  "
$ws.Range("E9").Value = "  +18.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.620.24"
$ws.Range("E10").Value = "  +7.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  +7.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.74"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000188"
$ws.Range("E14").Value = "  +6.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.102.16"
$ws.Range("E15").Value = "  +7.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "73.500.13"
$ws.Range("E16").Value = "  +6.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.18"
$ws.Range("E17").Value = "  +12.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.622.09"
$ws.Range("E18").Value = "  +7.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.04"
$ws.Range("E19").Value = "  +29.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  +11.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.50"
$ws.Range("E21").Value = "  +9.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("E22").Value = "  +18.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("E23").Value = "  +6.40%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.77"
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.12"
$ws.Range("E26").Value = "  +11.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.36"
$ws.Range("E27").Value = "  +13.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.758.83"
$ws.Range("E28").Value = "  +7.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0942"
$ws.Range("E30").Value = "  +14.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "521.43"
$ws.Range("E31").Value = "  +21.63%  "
$ws.Range("E32").Value = "  +19.58%  "
$ws.Range("E33").Value = "  +7.13%  "
$ws.Range("E34").Value = "  +8.52%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.13"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  +12.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.14"
$ws.Range("E38").Value = "  +6.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.24"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.89"
$ws.Range("E41").Value = "  +12.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.329"
$ws.Range("E42").Value = "  +9.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.66"
$ws.Range("E43").Value = "  +10.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "160.79"
$ws.Range("E44").Value = "  +23.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("E45").Value = "  +10.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.36"
$ws.Range("E46").Value = "  +14.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.94"
$ws.Range("E47").Value = "  +3.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0851"
$ws.Range("E48").Value = "  +18.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.62"
$ws.Range("E49").Value = "  +8.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.527"
$ws.Range("E50").Value = "  +9.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.53"
$ws.Range("E51").Value = "  +21.39%  "
